$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select row 169 (entire row) and delete it, shifting rows 170:223 up by one.
$ws.Rows.Item(169).Select()
$ws.Rows.Item(169).Delete()

# Restore selection/view state to match the post-edit workbook:
# topLeftCell="H148", selection activeCell="H169" sqref="A169:XFD169"
$ws.Range("H169").Activate()
$ws.Range("A169:XFD169").Select()
$excel.ActiveWindow.ScrollRow = 148
$excel.ActiveWindow.ScrollColumn = 8
